$wb = $excel.ActiveWorkbook

# --- Sheet "Country": update A2, remove B2 ---
$wsCountry = $wb.Worksheets.Item("Country")
$wsCountry.Range("A2").Value = 31.55
$wsCountry.Range("B2").ClearContents()

# --- Sheet "States": re-sorted rows 2-26 (descending CONCERN.HIGH), drop column D (MoM) ---
$wsStates = $wb.Worksheets.Item("States")
$wsStates.Cells.Item(2,1).Value = "IN-BR"
$wsStates.Cells.Item(2,2).Value = "Bihar"
$wsStates.Cells.Item(2,3).Value = 89.47
$wsStates.Cells.Item(3,1).Value = "IN-JH"
$wsStates.Cells.Item(3,2).Value = "Jharkhand"
$wsStates.Cells.Item(3,3).Value = 75
$wsStates.Cells.Item(4,1).Value = "IN-PB"
$wsStates.Cells.Item(4,2).Value = "Punjab"
$wsStates.Cells.Item(4,3).Value = 63.64
$wsStates.Cells.Item(5,1).Value = "IN-ML"
$wsStates.Cells.Item(5,2).Value = "Meghalaya"
$wsStates.Cells.Item(5,3).Value = 63.64
$wsStates.Cells.Item(6,1).Value = "IN-UP"
$wsStates.Cells.Item(6,2).Value = "Uttar Pradesh"
$wsStates.Cells.Item(6,3).Value = 62.67
$wsStates.Cells.Item(7,1).Value = "IN-MN"
$wsStates.Cells.Item(7,2).Value = "Manipur"
$wsStates.Cells.Item(7,3).Value = 56.25
$wsStates.Cells.Item(8,1).Value = "IN-MZ"
$wsStates.Cells.Item(8,2).Value = "Mizoram"
$wsStates.Cells.Item(8,3).Value = 54.55
$wsStates.Cells.Item(9,1).Value = "IN-DD"
$wsStates.Cells.Item(9,2).Value = "Daman and Diu"
$wsStates.Cells.Item(9,3).Value = 50
$wsStates.Cells.Item(10,1).Value = "IN-TR"
$wsStates.Cells.Item(10,2).Value = "Tripura"
$wsStates.Cells.Item(10,3).Value = 50
$wsStates.Cells.Item(11,1).Value = "IN-PY"
$wsStates.Cells.Item(11,2).Value = "Puducherry"
$wsStates.Cells.Item(11,3).Value = 50
$wsStates.Cells.Item(12,1).Value = "IN-DL"
$wsStates.Cells.Item(12,2).Value = "Delhi"
$wsStates.Cells.Item(12,3).Value = 45.45
$wsStates.Cells.Item(13,1).Value = "IN-AR"
$wsStates.Cells.Item(13,2).Value = "Arunachal Pradesh"
$wsStates.Cells.Item(13,3).Value = 44
$wsStates.Cells.Item(14,1).Value = "IN-NL"
$wsStates.Cells.Item(14,2).Value = "Nagaland"
$wsStates.Cells.Item(14,3).Value = 36.36
$wsStates.Cells.Item(15,1).Value = "IN-TS"
$wsStates.Cells.Item(15,2).Value = "Telangana"
$wsStates.Cells.Item(15,3).Value = 33.33
$wsStates.Cells.Item(16,1).Value = "IN-MP"
$wsStates.Cells.Item(16,2).Value = "Madhya Pradesh"
$wsStates.Cells.Item(16,3).Value = 30.77
$wsStates.Cells.Item(17,1).Value = "IN-OR"
$wsStates.Cells.Item(17,2).Value = "Odisha"
$wsStates.Cells.Item(17,3).Value = 30
$wsStates.Cells.Item(18,1).Value = "IN-JK"
$wsStates.Cells.Item(18,2).Value = "Jammu and Kashmir"
$wsStates.Cells.Item(18,3).Value = 27.27
$wsStates.Cells.Item(19,1).Value = "IN-HR"
$wsStates.Cells.Item(19,2).Value = "Haryana"
$wsStates.Cells.Item(19,3).Value = 27.27
$wsStates.Cells.Item(20,1).Value = "IN-AS"
$wsStates.Cells.Item(20,2).Value = "Assam"
$wsStates.Cells.Item(20,3).Value = 24.24
$wsStates.Cells.Item(21,1).Value = "IN-RJ"
$wsStates.Cells.Item(21,2).Value = "Rajasthan"
$wsStates.Cells.Item(21,3).Value = 15.15
$wsStates.Cells.Item(22,1).Value = "IN-WB"
$wsStates.Cells.Item(22,2).Value = "West Bengal"
$wsStates.Cells.Item(22,3).Value = 13.04
$wsStates.Cells.Item(23,1).Value = "IN-HP"
$wsStates.Cells.Item(23,2).Value = "Himachal Pradesh"
$wsStates.Cells.Item(23,3).Value = 8.33
$wsStates.Cells.Item(24,1).Value = "IN-CT"
$wsStates.Cells.Item(24,2).Value = "Chhattisgarh"
$wsStates.Cells.Item(24,3).Value = 7.41
$wsStates.Cells.Item(25,1).Value = "IN-GJ"
$wsStates.Cells.Item(25,2).Value = "Gujarat"
$wsStates.Cells.Item(25,3).Value = 3.03
$wsStates.Cells.Item(26,1).Value = "IN-MH"
$wsStates.Cells.Item(26,2).Value = "Maharashtra"
$wsStates.Cells.Item(26,3).Value = 2.78
$wsStates.Range("D2:D26").ClearContents()

# --- Sheet "Dark clusters": update B3, remove C2 and C3 ---
$wsDark = $wb.Worksheets.Item("Dark clusters")
$wsDark.Range("C2").ClearContents()
$wsDark.Range("B3").Value = 42.61
$wsDark.Range("C3").ClearContents()
